$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Productdata")
$ws.Range("E2").Value = 11.7029088
$ws.Range("E3").Value = 2.1871904
$ws.Range("E4").Value = 0.891504
$ws.Range("E5").Value = 1.414512
$ws.Range("E6").Value = 0.9700991999999999
$ws.Range("E7").Value = 0.2920128
$ws.Range("E8").Value = 0.0988416
$ws.Range("C9").Value = 839
$ws.Range("E9").Value = 0.8602304
$ws.Range("C10").Value = 607
$ws.Range("E10").Value = 0.5215104
$ws.Range("C11").Value = 1803
$ws.Range("E11").Value = 0.757344
$ws.Range("C12").Value = 418
$ws.Range("E12").Value = 1.2192928
$ws.Range("C13").Value = 3237
$ws.Range("E13").Value = 12.58516799999999
$ws.Range("C14").Value = 1211
$ws.Range("E14").Value = 4.880228800000001
$ws.Range("C15").Value = 300
$ws.Range("E15").Value = 0.8702496
$ws.Range("C16").Value = 413
$ws.Range("E16").Value = 0.8702144000000001
$ws.Range("C17").Value = 609
$ws.Range("E17").Value = 1.415232
$ws.Range("C18").Value = 181
$ws.Range("E18").Value = 0.457632
$ws.Range("C19").Value = 60
$ws.Range("E19").Value = 0.1367808
$ws.Range("E20").Value = 63.17891200000001
$ws.Range("E21").Value = 67.0705408
$ws.Range("E22").Value = 83.02141440000001
$ws.Range("E23").Value = 255.4542016

$ws = $wb.Worksheets.Item("ForecastedAverageDemand")
$ws.Range("C2").Value = 594
$ws.Range("D2").Value = 137
$ws.Range("F2").Value = 300
$ws.Range("G2").Value = 93
$ws.Range("H2").Value = 30
$ws.Range("I2").Value = 418
$ws.Range("J2").Value = 302
$ws.Range("K2").Value = 895
$ws.Range("L2").Value = 206
$ws.Range("C3").Value = 603
$ws.Range("D3").Value = 150
$ws.Range("F3").Value = 297
$ws.Range("G3").Value = 90
$ws.Range("H3").Value = 40
$ws.Range("J3").Value = 305
$ws.Range("K3").Value = 908
$ws.Range("L3").Value = 212
$ws.Range("C4").Value = 610
$ws.Range("D4").Value = 145
$ws.Range("F4").Value = 305
$ws.Range("G4").Value = 91
$ws.Range("H4").Value = 25
$ws.Range("K4").Value = 903
$ws.Range("L4").Value = 211
$ws.Range("C5").Value = 601
$ws.Range("D5").Value = 155
$ws.Range("F5").Value = 304
$ws.Range("G5").Value = 90
$ws.Range("H5").Value = 35
$ws.Range("I5").Value = 415
$ws.Range("J5").Value = 300
$ws.Range("K5").Value = 903
$ws.Range("L5").Value = 202
$ws.Range("C6").Value = 603
$ws.Range("D6").Value = 151
$ws.Range("F6").Value = 306
$ws.Range("G6").Value = 90
$ws.Range("H6").Value = 26
$ws.Range("I6").Value = 420
$ws.Range("J6").Value = 304
$ws.Range("K6").Value = 899
$ws.Range("L6").Value = 203

$ws = $wb.Worksheets.Item("ForcastedStandardDeviation")
$ws.Range("C2").Value = 74.25
$ws.Range("D2").Value = 17.125
$ws.Range("F2").Value = 37.5
$ws.Range("G2").Value = 11.625
$ws.Range("H2").Value = 3.75
$ws.Range("I2").Value = 52.25
$ws.Range("J2").Value = 37.75
$ws.Range("K2").Value = 111.875
$ws.Range("L2").Value = 25.75
$ws.Range("C3").Value = 113.0625
$ws.Range("D3").Value = 28.125
$ws.Range("F3").Value = 55.6875
$ws.Range("G3").Value = 16.875
$ws.Range("H3").Value = 7.5
$ws.Range("J3").Value = 57.1875
$ws.Range("K3").Value = 170.25
$ws.Range("L3").Value = 39.75
$ws.Range("C4").Value = 133.4375
$ws.Range("D4").Value = 31.71875
$ws.Range("F4").Value = 66.71875
$ws.Range("G4").Value = 19.90625
$ws.Range("H4").Value = 5.46875
$ws.Range("K4").Value = 197.53125
$ws.Range("L4").Value = 46.15625
$ws.Range("C5").Value = 140.859375
$ws.Range("D5").Value = 36.328125
$ws.Range("F5").Value = 71.25
$ws.Range("G5").Value = 21.09375
$ws.Range("H5").Value = 8.203125
$ws.Range("I5").Value = 97.265625
$ws.Range("J5").Value = 70.3125
$ws.Range("K5").Value = 211.640625
$ws.Range("L5").Value = 47.34375
$ws.Range("C6").Value = 146.0390625
$ws.Range("D6").Value = 36.5703125
$ws.Range("F6").Value = 74.109375
$ws.Range("G6").Value = 21.796875
$ws.Range("H6").Value = 6.296875
$ws.Range("I6").Value = 101.71875
$ws.Range("J6").Value = 73.625
$ws.Range("K6").Value = 217.7265625
$ws.Range("L6").Value = 49.1640625

$ws = $wb.Worksheets.Item("Capacity")
$ws.Range("B2").Value = 20165.6
$ws.Range("B3").Value = 210196
$ws.Range("B4").Value = 210196
$ws.Range("B5").Value = 1981848
